# Update the "想去人数" (want-to-go count) figures in column F
# for the "展览" (Exhibition) sheet and the "全部类型" (All types)
# aggregate sheet, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new value for column F
$sheet1Updates = @{
    2  = 1208
    3  = 430
    4  = 292
    6  = 16
    7  = 12409
    8  = 64
    10 = 20
    12 = 176
    13 = 12247
    14 = 4854
    15 = 4739
    16 = 141
    20 = 955
    22 = 365
    23 = 170
    25 = 5216
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" - row => new value for column F
$sheet4Updates = @{
    2  = 1208
    3  = 430
    4  = 292
    8  = 16
    9  = 12409
    10 = 64
    12 = 20
    14 = 176
    15 = 12247
    16 = 4854
    17 = 4739
    18 = 141
    22 = 955
    24 = 365
    25 = 170
    27 = 5216
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
